# Auto-generated script applying diff changes to Alpha_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2645.25
$ws.Range("J17").Value = 2645.25
$ws.Range("L17").Value = 7935.75
$ws.Range("N17").Value = -8271.75

$ws.Range("H18").Value = 401
$ws.Range("I18").Value = 401
$ws.Range("K18").Value = 401
$ws.Range("M18").Value = -117

$ws.Range("H28").Value = 1500
$ws.Range("I28").Value = 1500
$ws.Range("K28").Value = 1500
$ws.Range("M28").Value = -1015

$ws.Range("H62").Value = 7799.8
$ws.Range("I62").Value = 7999.75
$ws.Range("K62").Value = 7999.75
$ws.Range("M62").Value = -7375.75

$ws.Range("H65").Value = 7799.8
$ws.Range("I65").Value = 7999.75
$ws.Range("K65").Value = 39998.75
$ws.Range("M65").Value = -36878.75

$ws.Range("H76").Value = 5333
$ws.Range("J76").Value = 5666
$ws.Range("L76").Value = 5666
$ws.Range("N76").Value = -6296

$ws.Range("H79").Value = 5333
$ws.Range("J79").Value = 5666
$ws.Range("L79").Value = 5666
$ws.Range("N79").Value = -7850

$ws.Range("H116").Value = 3385.7144
$ws.Range("I116").Value = 3350
$ws.Range("K116").Value = 3350
$ws.Range("M116").Value = 92

$ws.Range("H129").Value = 1980.2041
$ws.Range("I129").Value = 464.14285
$ws.Range("J129").Value = 2096.8242
$ws.Range("K129").Value = 1392.42855
$ws.Range("L129").Value = 6290.4726
$ws.Range("M129").Value = 3607.57145
$ws.Range("N129").Value = -16290.4726

$ws.Range("H137").Value = 2124.8215
$ws.Range("I137").Value = 1849.1765
$ws.Range("K137").Value = 5547.529500000001
$ws.Range("M137").Value = -2997.529500000001

$ws.Range("H138").Value = 2478.875
$ws.Range("I138").Value = 1121
$ws.Range("J138").Value = 3836.75
$ws.Range("K138").Value = 3363
$ws.Range("L138").Value = 11510.25
$ws.Range("M138").Value = 1777
$ws.Range("N138").Value = -21790.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4222.5137
$ws.Range("I32").Value = 4222.5137
$ws.Range("K32").Value = 4222.5137
$ws.Range("M32").Value = -3935.5137

$ws.Range("H45").Value = 1657.2
$ws.Range("J45").Value = 1999
$ws.Range("L45").Value = 1999
$ws.Range("N45").Value = -2753

$ws.Range("H74").Value = 2356.5715
$ws.Range("I74").Value = 1994.2941
$ws.Range("K74").Value = 1994.2941
$ws.Range("M74").Value = -1120.2941

$ws.Range("H77").Value = 2356.5715
$ws.Range("I77").Value = 1994.2941
$ws.Range("K77").Value = 9971.470499999999
$ws.Range("M77").Value = -5603.470499999999

$ws.Range("H109").Value = 49950
$ws.Range("J109").Value = 49950
$ws.Range("L109").Value = 49950
$ws.Range("N109").Value = -52724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1642.7222
$ws.Range("I99").Value = 1444.5883
$ws.Range("J99").Value = 5011
$ws.Range("K99").Value = 1444.5883
$ws.Range("L99").Value = 5011
$ws.Range("M99").Value = 53.41170000000011
$ws.Range("N99").Value = -8007

$ws.Range("H134").Value = 4589.24
$ws.Range("I134").Value = 4562.2607
$ws.Range("K134").Value = 13686.7821
$ws.Range("M134").Value = -11151.7821

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").ClearContents()
$ws.Range("N55").Value = 0

$ws.Range("H134").Value = 4656.88
$ws.Range("I134").Value = 3233.3684
$ws.Range("K134").Value = 9700.1052
$ws.Range("M134").Value = -7165.1052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 603.34784
$ws.Range("J113").Value = 612.63635
$ws.Range("L113").Value = 1837.90905
$ws.Range("N113").Value = -6177.90905

$ws.Range("H121").Value = 2361269
$ws.Range("I121").Value = 113060.445
$ws.Range("K121").Value = 339181.335
$ws.Range("M121").Value = -337871.335

$ws.Range("H122").Value = 1973.625
$ws.Range("I122").Value = 1630
$ws.Range("J122").Value = 2179.8
$ws.Range("K122").Value = 14670
$ws.Range("L122").Value = 19618.2
$ws.Range("M122").Value = -12220
$ws.Range("N122").Value = -24518.2

$ws.Range("H128").Value = 138155.33
$ws.Range("I128").Value = 138155.33
$ws.Range("K128").Value = 414465.99
$ws.Range("M128").Value = -409485.99

$ws.Range("H137").Value = 6059.4614
$ws.Range("J137").Value = 5987.4
$ws.Range("L137").Value = 17962.2
$ws.Range("N137").Value = -28162.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 136.53847
$ws.Range("I2").Value = 127.083336
$ws.Range("K2").Value = 127.083336
$ws.Range("M2").Value = -14.083336

$ws.Range("H43").Value = 7815
$ws.Range("I43").Value = 2450.8333
$ws.Range("K43").Value = 2450.8333
$ws.Range("M43").Value = -2299.8333

$ws.Range("H102").Value = 1176.5714
$ws.Range("J102").Value = 1144
$ws.Range("L102").Value = 1144
$ws.Range("N102").Value = -4388

$ws.Range("H132").Value = 2924.9092
$ws.Range("I132").Value = 2782.1428
$ws.Range("J132").Value = 3174.75
$ws.Range("K132").Value = 8346.428400000001
$ws.Range("L132").Value = 9524.25
$ws.Range("M132").Value = -5816.428400000001
$ws.Range("N132").Value = -14584.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1033.7858
$ws.Range("I16").Value = 1096.3846
$ws.Range("J16").Value = 220
$ws.Range("K16").Value = 1096.3846
$ws.Range("L16").Value = 220
$ws.Range("M16").Value = -926.3846000000001
$ws.Range("N16").Value = -560

$ws.Range("H46").Value = 3216.6667
$ws.Range("I46").Value = 703.5714
$ws.Range("K46").Value = 703.5714
$ws.Range("M46").Value = -515.5714

$ws.Range("H68").Value = 4000
$ws.Range("I68").Value = 4000
$ws.Range("K68").Value = 4000
$ws.Range("M68").Value = -3251

$ws.Range("H71").Value = 4000
$ws.Range("I71").Value = 4000
$ws.Range("K71").Value = 20000
$ws.Range("M71").Value = -16256

$ws.Range("H132").Value = 3726.3333
$ws.Range("I132").Value = 2824.8
$ws.Range("K132").Value = 8474.400000000001
$ws.Range("M132").Value = -5944.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 8000
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H40").Value = 13496.75
$ws.Range("J40").Value = 13496.75
$ws.Range("L40").Value = 13496.75
$ws.Range("N40").Value = -13794.75

$ws.Range("H42").Value = 49999
$ws.Range("I42").Value = 49999
$ws.Range("K42").Value = 49999
$ws.Range("M42").Value = -49621

$ws.Range("H81").Value = 2515.5789
$ws.Range("I81").Value = 2239.7
$ws.Range("J81").Value = 2822.111
$ws.Range("K81").Value = 4479.4
$ws.Range("L81").Value = 5644.222
$ws.Range("M81").Value = -3418.4
$ws.Range("N81").Value = -7766.222

$ws.Range("H84").Value = 2515.5789
$ws.Range("I84").Value = 2239.7
$ws.Range("J84").Value = 2822.111
$ws.Range("K84").Value = 22397
$ws.Range("L84").Value = 28221.11
$ws.Range("M84").Value = -17093
$ws.Range("N84").Value = -38829.11

$ws.Range("H126").Value = 5220.8
$ws.Range("I126").Value = 8250
$ws.Range("K126").Value = 24750
$ws.Range("M126").Value = -22280

$ws.Range("H136").Value = 5336
$ws.Range("I136").Value = 5336
$ws.Range("K136").Value = 16008
$ws.Range("M136").Value = -13458
